$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Cells.Item(1, 8)
Write-Host ("Cells.Item(1,8) Value2: " + $c.Value2)
Write-Host ("Address: " + $c.Address())
